$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column B values first (in row order), then column C values,
# matching the order new shared strings were introduced.
$ws.Range("B2").Value = "vägen 2, adress continued."
$ws.Range("B11").Value = "mera adresss"
$ws.Range("B13").Value = "här bor jag"
$ws.Range("B30").Value = "och jag bor här"

$ws.Range("C15").Value = "skaldjur"
$ws.Range("C19").Value = "fisk"
$ws.Range("C24").Value = "köttbullar"
$ws.Range("C3").Value = "ingefära"

# Size column B to fit its longest new entry (matches the bestFit width
# Excel computed for "vagen 2, adress continued." in the authored file).
$ws.Columns.Item(2).ColumnWidth = 24

# Set the active selection to C3 as in the edited file
$ws.Range("C3").Select() | Out-Null
